$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.000457525253296
$ws.Range("B1").Value = 1.324018478393555
$ws.Range("C1").Value = 5.674402236938477
$ws.Range("D1").Value = 1.667941570281982
$ws.Range("E1").Value = 1.022565245628357
